$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$newRow.Cells(1).Range.Text = "7"
$newRow.Cells(2).Range.Text = "-"
$newRow.Cells(3).Range.Text = "0.377"
